$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.709.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = "'" + '1.879.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.20%  '
$ws.Range("D4").Value = "'" + '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'" + '237.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.88%  '
$ws.Range("D6").Value = "'" + '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = "'" + '0.4750'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.30%  '
$ws.Range("D8").Value = "'" + '0.2830'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.90%  '
$ws.Range("D9").Value = "'" + '0.06517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.99%  '
$ws.Range("D10").Value = "'" + '18.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +14.73%  '
$ws.Range("D11").Value = "'" + '1.872.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("D12").Value = "'" + '0.07584'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("D13").Value = "'" + '95.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +14.72%  '
$ws.Range("D14").Value = "'" + '5.061'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").Value = "'" + '0.6496'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.99%  '
$ws.Range("D16").Value = "'" + '303.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +33.44%  '
$ws.Range("D17").Value = "'" + '30.682.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").Value = "'" + '0.9999'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = "'" + '13.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.17%  '
$ws.Range("D20").Value = "'" + '0.000007522'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.29%  '
$ws.Range("D21").Value = "'" + '2.145.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").Value = "'" + '0.9991'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = "'" + '5.137'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.09%  '
$ws.Range("D24").Value = "'" + '6.152'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.09%  '
$ws.Range("D25").Value = "'" + '169.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.34%  '
$ws.Range("D26").Value = "'" + '9.254'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").Value = "'" + '19.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.87%  '
$ws.Range("D28").Value = "'" + '1.952'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.71%  '
$ws.Range("D29").Value = "'" + '0.1058'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("D30").Value = "'" + '1.364'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").Value = "'" + '4.155'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").Value = "'" + '3.953'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.09%  '
$ws.Range("D33").Value = "'" + '0.05011'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.67%  '
$ws.Range("D34").Value = "'" + '1.170'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.98%  '
$ws.Range("D35").Value = "'" + '0.7195'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").Value = "'" + '2.707'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("D37").Value = "'" + '0.01914'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("D38").Value = "'" + '2.695'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = "'" + '2.048'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.18%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'" + '0.8989'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").Value = "'" + '107.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").Value = "'" + '1.0000'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").Value = "'" + '0.4188'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.60%  '
$ws.Range("D44").Value = "'" + '5.587'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("D45").Value = "'" + '7.329'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.28%  '
$ws.Range("D46").Value = "'" + '65.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.66%  '
$ws.Range("D47").Value = "'" + '8.983'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.68%  '
$ws.Range("D48").Value = "'" + '0.1214'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").Value = "'" + '34.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.22%  '
$ws.Range("D50").Value = "'" + '0.05584'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("E51").Value = '  +2.42%  '
